$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 10: Objetivos: value changes from the long description to the
# professor reference text.
# ---------------------------------------------------------------------------
$ws.Cells.Item(10, 2).Value = "144651 - Antonio Fernando Sartori"
$ws.Cells.Item(10, 3).Value = "144651 - Antonio Fernando Sartori"

# ---------------------------------------------------------------------------
# Row 13 becomes "Programa resumido:" with the activation-date text reused
# in B/C (copied from B8/C8 so the cell keeps its original text typing
# instead of Excel re-interpreting "01/01/2020" as a date value).
# ---------------------------------------------------------------------------
$ws.Cells.Item(13, 1).Value = "Programa resumido:"
$ws.Rows.Item(13).RowHeight = 60

$ws.Cells.Item(8, 2).Copy()
$ws.Cells.Item(13, 2).PasteSpecial(-4163)
$ws.Cells.Item(8, 3).Copy()
$ws.Cells.Item(13, 3).PasteSpecial(-4163)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Row 14 becomes "Short syllabus:" with no B/C value.
# ---------------------------------------------------------------------------
$ws.Cells.Item(14, 1).Value = "Short syllabus:"
$ws.Cells.Item(14, 2).ClearContents()
$ws.Cells.Item(14, 3).ClearContents()
$ws.Rows.Item(14).RowHeight = 60

# ---------------------------------------------------------------------------
# Row 15 becomes "Programa:" with the professor reference text reused again.
# ---------------------------------------------------------------------------
$ws.Cells.Item(15, 1).Value = "Programa:"
$ws.Cells.Item(15, 2).Value = "144651 - Antonio Fernando Sartori"
$ws.Cells.Item(15, 3).Value = "144651 - Antonio Fernando Sartori"
$ws.Rows.Item(15).RowHeight = 120

# ---------------------------------------------------------------------------
# Row 16 becomes "Syllabus:" with no B/C value.
# ---------------------------------------------------------------------------
$ws.Cells.Item(16, 1).Value = "Syllabus:"
$ws.Cells.Item(16, 2).ClearContents()
$ws.Cells.Item(16, 3).ClearContents()
$ws.Rows.Item(16).RowHeight = 120

# ---------------------------------------------------------------------------
# Row 17 becomes "Avaliação:" with no B/C value and default row height.
# ---------------------------------------------------------------------------
$ws.Cells.Item(17, 1).Value = "Avaliação:"
$ws.Cells.Item(17, 2).ClearContents()
$ws.Cells.Item(17, 3).ClearContents()
$ws.Rows.Item(17).RowHeight = 15

# ---------------------------------------------------------------------------
# Row 18 becomes "Método:" with the "3577649 - Carlos Angelo Nunes" text.
# ---------------------------------------------------------------------------
$ws.Cells.Item(18, 1).Value = "Método:"
$ws.Cells.Item(18, 2).Value = "3577649 - Carlos Angelo Nunes"
$ws.Cells.Item(18, 3).Value = "3577649 - Carlos Angelo Nunes"
$ws.Rows.Item(18).RowHeight = 60

# ---------------------------------------------------------------------------
# Row 19 becomes "Critério:" with the grading-method paragraph (previously
# on row 20).
# ---------------------------------------------------------------------------
$ws.Cells.Item(19, 1).Value = "Critério:"
$ws.Cells.Item(19, 2).Value = "Serão usadas duas notas P1 e P2. A P1 será uma prova escrita e a P2 será a soma de uma nota de seminário e uma nota de trabalho escrito."
$ws.Cells.Item(19, 3).Value = "Serão usadas duas notas P1 e P2. A P1 será uma prova escrita e a P2 será a soma de uma nota de seminário e uma nota de trabalho escrito."
$ws.Rows.Item(19).RowHeight = 60

# ---------------------------------------------------------------------------
# Row 20 becomes "Norma de recuperação:" with the grade-criteria paragraph
# (previously on row 21).
# ---------------------------------------------------------------------------
$ws.Cells.Item(20, 1).Value = "Norma de recuperação:"
$ws.Cells.Item(20, 2).Value = "Média Final(MF) = (P1 + 2P2)/3 MF menor que 3,0 - reprovado. MF maior,ou igual, a 3,0 até menor que 5,0 - recuperação. MF maior, ou igual, a 5,0 - aprovado."
$ws.Cells.Item(20, 3).Value = "Média Final(MF) = (P1 + 2P2)/3 MF menor que 3,0 - reprovado. MF maior,ou igual, a 3,0 até menor que 5,0 - recuperação. MF maior, ou igual, a 5,0 - aprovado."
$ws.Rows.Item(20).RowHeight = 60

# ---------------------------------------------------------------------------
# Row 21 becomes "Bibliografia:" with the recovery-norm paragraph
# (previously on row 22), and a taller row height.
# ---------------------------------------------------------------------------
$ws.Cells.Item(21, 1).Value = "Bibliografia:"
$ws.Cells.Item(21, 2).Value = "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2."
$ws.Cells.Item(21, 3).Value = "Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2."
$ws.Rows.Item(21).RowHeight = 120

# ---------------------------------------------------------------------------
# Row 22 becomes "Requisitos:" with no B/C value and default row height
# (the old Bibliografia long text is dropped entirely).
# ---------------------------------------------------------------------------
$ws.Cells.Item(22, 1).Value = "Requisitos:"
$ws.Cells.Item(22, 2).ClearContents()
$ws.Cells.Item(22, 3).ClearContents()
$ws.Rows.Item(22).RowHeight = 15

# ---------------------------------------------------------------------------
# Row 23 keeps only the requirement text in B/C (previously on row 25),
# with column A cleared, and the old rows 24/25 are removed entirely.
# ---------------------------------------------------------------------------
$ws.Cells.Item(23, 1).ClearContents()
$ws.Cells.Item(23, 2).Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)`n"
$ws.Cells.Item(23, 3).Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)`n"
$ws.Rows.Item(23).RowHeight = 30

$ws.Rows.Item(25).Delete()
$ws.Rows.Item(24).Delete()
